$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.439.71'
$ws.Range('E2').Value = '  -3.02%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.418.83'
$ws.Range('E3').Value = '  +6.94%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '294.60'
$ws.Range('E5').Value = '  -2.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.52'
$ws.Range('E6').Value = '  -6.09%  '
$ws.Range('E7').Value = '  +0.91%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.502'
$ws.Range('E9').Value = '  -0.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.96'
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0780'
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.04'
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.789.88'
$ws.Range('E14').Value = '  +7.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.421.29'
$ws.Range('E15').Value = '  +7.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.836'
$ws.Range('E16').Value = '  +5.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.08'
$ws.Range('E17').Value = '  +3.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '45.349.10'
$ws.Range('E18').Value = '  -3.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.36'
$ws.Range('E19').Value = '  -2.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0940'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.15'
$ws.Range('E21').Value = '  +6.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.05'
$ws.Range('E22').Value = '  +2.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '241.08'
$ws.Range('E23').Value = '  -2.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.78'
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('E26').Value = '  +3.66%  '
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.35'
$ws.Range('E28').Value = '  -9.03%  '
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '21.20'
$ws.Range('E30').Value = '  +5.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.75'
$ws.Range('E31').Value = '  +15.97%  '
$ws.Range('E32').Value = '  -2.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '148.37'
$ws.Range('E33').Value = '  +1.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.45'
$ws.Range('E34').Value = '  +1.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0763'
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.113'
$ws.Range('E36').Value = '  -1.28%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.97'
$ws.Range('E37').Value = '  +17.27%  '
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.72'
$ws.Range('E39').Value = '  -8.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.78'
$ws.Range('E40').Value = '  -2.37%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.009.56'
$ws.Range('E42').Value = '  +13.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.22'
$ws.Range('E43').Value = '  +3.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '88.96'
$ws.Range('E45').Value = '  -2.69%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.75'
$ws.Range('E46').Value = '  -9.11%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.96'
$ws.Range('E47').Value = '  +24.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.56'
$ws.Range('E48').Value = '  +9.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '100.77'
$ws.Range('E49').Value = '  +7.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.657.97'
$ws.Range('E50').Value = '  +6.96%  '
$ws.Range('E51').Value = '  -1.05%  '
